$wb = $excel.ActiveWorkbook

# Final patient/measure identifiers (column A) and their hash values (column B)
# for rows 2-21, identical on both worksheets.
$data = @(
    @("AMI1,DEN_0142_0164_0639", "caf50e70e548d61d097c68e9001ded60"),
    @("AMI1,NUM_0142_0164_0639", "649910f4a9da359b17d7c1b012decc17"),
    @("AMI2,DEN_0163", "b03719c7f99502b7990918baf4640f70"),
    @("AMI2,NUM_0163", "116c5a883ccdc89b2531bcbae3a403ab"),
    @("ASTH,DEN_0338_496", "3f3a63ed0f7c69e7d39e0f7f6f123151"),
    @("ASTH,NUM_0338", "29778d9cb06701de524b752c1f4de992"),
    @("CAP,DEN_0147_372", "4f2694161da81fb6d16f6fb421239e3d"),
    @("CAP,NUM_0147_372", "e9f710b1f129c10bf69d3c29432e166f"),
    @("NEWB,DEN_0480_0716_1354", "1ef57c0a5ac2ef5a2b50b3f4bb04d76c"),
    @("NEWB,NUM_0480_0716_1354", "f03f9bf5cc4bfa7298ae1a4804cc7e5f"),
    @("OB,DEN_0469_376", "52e4994176da107bf82584a316146e68"),
    @("OB,NUM_0469_0376", "f8745b14e07f41c49d8785a1d9920b77"),
    @("SCIP,DEN_0453_0527_0528", "4419c3df96d1e9699136eaa29b9cd44b"),
    @("SCIP,NUM_0453_0527_0528", "1d498bdc11d16adefd411599c04be5ee"),
    @("STROKE,DEN_0435_0436_0437_ 0438_ 0439_ 0440_ 0441", "879223b7a0d9fee24b76afac7e1ce268"),
    @("STROKE,NUM_0435_0436_0437_ 0439_ 0440_ 0441", "e8932b4783869759f5217d12cc6efd19"),
    @("STROKE,NUM_ 0438", "cce692bcfd3d4dc602a052190bb30d8b"),
    @("VTE1,NUM_0371", "59074d1606548365d4db75931397082c"),
    @("VTE2,DEN_0373_0374_0375", "c6bbb7342c4cf4ebeeaf1a417646db69"),
    @("VTE2,NUM_0373_0374_0375", "6d4da8955b5aee38deb098909d3c4a06")
)

for ($s = 1; $s -le 2; $s++) {
    $ws = $wb.Worksheets.Item($s)

    for ($i = 0; $i -lt $data.Length; $i++) {
        $row = $i + 2
        $ws.Cells.Item($row, 1).Value = $data[$i][0]
        $ws.Cells.Item($row, 2).Value = $data[$i][1]
    }

    # Column A picks up the same style already used by column B (numFmtId 1,
    # linked to the Normal cell style) for the whole data block.
    $ws.Range("B2:B21").Copy()
    $ws.Range("A2:A21").PasteSpecial(-4122)
    $excel.CutCopyMode = $false

    $ws.Range("A2:B21").Select()
}

$ws1 = $wb.Worksheets.Item(1)
$ws1.Activate()
$ws1.Range("A2:B21").Select()

$ws2 = $wb.Worksheets.Item(2)
$ws2.Activate()
$ws2.Range("A2:B21").Select()
